# Reformat the shared-string payload in A2 into pretty-printed JSON (with
# \uXXXX escapes for non-ASCII punctuation) and move it up into A1, then
# drop the old numeric placeholder row (A1=0) that used to carry the
# bold+bordered header style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$text = 'questions = [
    {
        "title": "Your team cannot make outbound calls using talk, the line cuts after the first ring. You have contacted the support team and they have explained that they will need access to your instance. How should you grant them access?",
        "ques_type": 2,
        "options": [
            "Account &gt Security &gt Advanced &gt  Account Assumption &gt enable",
            "People &gt Edit email address &gt Change email address with Zendesk customer support agent email address. ",
            "User\u2019s profile &gt Adopt Identity",
            "Account &gt Security &gt Advanced &gt Authentication &gt Require two-factor authentication. "
        ],
        "score": "Account &gt Security &gt Advanced &gt  Account Assumption &gt enable"
    },
    {
        "title": "Your website displayed the wrong price for one of your products for a few days last week. Many clients bought the product during that time, and have since contacted you to ask about the delivery. You need to let them know that the transaction will be canceled and they will be refunded. How should you reply to all of them to solve the issue?",
        "ques_type": 2,
        "options": [
            "Add other team members as followers so they can help you to provide replies to the customers. ",
            "Use a shortcut to find the tickets faster.",
            "Bulk edit the tickets to apply the same answers.",
            "Merge all the tickets using the merge option."
        ],
        "score": "Bulk edit the tickets to apply the same answers."
    },
    {
        "title": "You have created a trigger to notify your customers that your office is closed for a national holiday. However, your clients are receiving two notifications: the default \u201cNotify UserRequest Received\u201d notification, and the new one you have created.Which action should you take to avoid sending the two notifications to end-users?",
        "ques_type": 2,
        "options": [
            "Delete the trigger by default \u201cnotify user request received\u201d.",
            "Delete the trigger \u201cnotify user holiday period\u201d you have created. ",
            "Update the trigger by default \u201cnotify user request received\u201d with \u201cday off // is // yes\u201c.",
            "Update the default trigger \u201cnotify user request received with \u201cday off // is // no\u201d."
        ],
        "score": "Update the default trigger \u201cnotify user request received with \u201cday off // is // no\u201d."
    },
    {
        "title": "You wish to compare the results of each team member in the support team. In order to, you need to gather information about the number of tickets handled by each agent in the previous month How should you gather this information?",
        "ques_type": 2,
        "options": [
            "Create an export from Support.",
            "Check each agent\u2019s profile activity.",
            "Create a shared view.",
            "Generate a report on Explore."
        ],
        "score": "Generate a report on Explore."
    }
]'

# Remove row 2 (the old shared-string cell) and shift the remaining cells
# up, so the text ends up alone on row 1.
$ws.Range("A2").Delete()

# Drop the bold/border/center-top style that used to live on A1 (the old
# numeric placeholder) before writing the new value into it.
$ws.Range("A1").ClearFormats()
$ws.Range("A1").Value = $text
